# Auto-generated Excel COM-interop edit script
# Applies the Slovakia 2-liga 2023-2024 data refresh:
# - 29 existing rows get their match/odds columns (F:V) replaced
#   with refreshed data (teams/scores/odds/timestamps/url reshuffled
#   among rows sharing the same match date).
# - 4 new rows (113-116) are appended for matches played 03-04/11/2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing rows (F:V) ----
# Row 3
$ws.Range("F3").Value = 'FK Humenne'
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 'Komarno'
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 2.83
$ws.Range("K3").Value = '28/07/2023 03:42'
$ws.Range("L3").Value = 2.97
$ws.Range("M3").Value = '29/07/2023 16:11'
$ws.Range("N3").Value = 3.09
$ws.Range("O3").Value = '28/07/2023 03:42'
$ws.Range("P3").Value = 3.19
$ws.Range("Q3").Value = '29/07/2023 16:11'
$ws.Range("R3").Value = 2.28
$ws.Range("S3").Value = '28/07/2023 03:42'
$ws.Range("T3").Value = 2.37
$ws.Range("U3").Value = '29/07/2023 16:11'
$ws.Range("V3").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/fk-humenne-komarno/dfEiFXZ8/'

# Row 5
$ws.Range("F5").Value = 'Spisska Nova Ves'
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 'Presov'
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 4.59
$ws.Range("K5").Value = '28/07/2023 03:42'
$ws.Range("L5").Value = 5.61
$ws.Range("M5").Value = '29/07/2023 15:53'
$ws.Range("N5").Value = 3.99
$ws.Range("O5").Value = '28/07/2023 03:42'
$ws.Range("P5").Value = 4.81
$ws.Range("Q5").Value = '29/07/2023 15:53'
$ws.Range("R5").Value = 1.57
$ws.Range("S5").Value = '28/07/2023 03:42'
$ws.Range("T5").Value = 1.46
$ws.Range("U5").Value = '29/07/2023 15:53'
$ws.Range("V5").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/spisska-nova-ves-presov/tnW4iUs2/'

# Row 7
$ws.Range("F7").Value = 'Petrzalka'
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 'Pohronie'
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2.3
$ws.Range("K7").Value = '28/07/2023 05:13'
$ws.Range("L7").Value = 2.42
$ws.Range("M7").Value = '28/07/2023 20:04'
$ws.Range("N7").Value = 3.3
$ws.Range("O7").Value = '28/07/2023 05:13'
$ws.Range("P7").Value = 3.37
$ws.Range("Q7").Value = '29/07/2023 15:03'
$ws.Range("R7").Value = 2.73
$ws.Range("S7").Value = '28/07/2023 05:13'
$ws.Range("T7").Value = 2.72
$ws.Range("U7").Value = '29/07/2023 09:34'
$ws.Range("V7").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/petrzalka-fk-pohronie/O8UXf3lf/'

# Row 8
$ws.Range("F8").Value = 'Povazska Bystrica'
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 'Samorin'
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 2.25
$ws.Range("K8").Value = '28/07/2023 03:42'
$ws.Range("L8").Value = 2.3
$ws.Range("M8").Value = '29/07/2023 16:24'
$ws.Range("N8").Value = 3.33
$ws.Range("O8").Value = '28/07/2023 03:42'
$ws.Range("P8").Value = 3.14
$ws.Range("Q8").Value = '29/07/2023 16:24'
$ws.Range("R8").Value = 2.71
$ws.Range("S8").Value = '28/07/2023 03:42'
$ws.Range("T8").Value = 3.16
$ws.Range("U8").Value = '29/07/2023 16:24'
$ws.Range("V8").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/povazska-bystrica-samorin/pIVTeqZm/'

# Row 18
$ws.Range("F18").Value = 'Puchov'
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 'Povazska Bystrica'
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 1.95
$ws.Range("K18").Value = '11/08/2023 05:13'
$ws.Range("L18").Value = 1.76
$ws.Range("M18").Value = '12/08/2023 16:51'
$ws.Range("N18").Value = 3.66
$ws.Range("O18").Value = '11/08/2023 05:13'
$ws.Range("P18").Value = 4.31
$ws.Range("Q18").Value = '12/08/2023 16:51'
$ws.Range("R18").Value = 3.12
$ws.Range("S18").Value = '11/08/2023 05:13'
$ws.Range("T18").Value = 3.67
$ws.Range("U18").Value = '12/08/2023 16:46'
$ws.Range("V18").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-povazska-bystrica/fDOESsAP/'

# Row 21
$ws.Range("F21").Value = 'Komarno'
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 'L. Mikulas'
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = 2.16
$ws.Range("K21").Value = '12/08/2023 08:42'
$ws.Range("L21").Value = 1.98
$ws.Range("M21").Value = '12/08/2023 16:59'
$ws.Range("N21").Value = 3.31
$ws.Range("O21").Value = '12/08/2023 08:42'
$ws.Range("P21").Value = 3.43
$ws.Range("Q21").Value = '12/08/2023 16:59'
$ws.Range("R21").Value = 3.1
$ws.Range("S21").Value = '12/08/2023 08:42'
$ws.Range("T21").Value = 3.62
$ws.Range("U21").Value = '12/08/2023 16:59'
$ws.Range("V21").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/komarno-l-mikulas/OShF7tYa/'

# Row 23
$ws.Range("F23").Value = 'Spisska Nova Ves'
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 'Zilina B'
$ws.Range("I23").Value = 2
$ws.Range("J23").Value = 3.05
$ws.Range("K23").Value = '12/08/2023 08:42'
$ws.Range("L23").Value = 2.7
$ws.Range("M23").Value = '12/08/2023 16:43'
$ws.Range("N23").Value = 3.81
$ws.Range("O23").Value = '12/08/2023 08:42'
$ws.Range("P23").Value = 3.86
$ws.Range("Q23").Value = '12/08/2023 16:27'
$ws.Range("R23").Value = 2.01
$ws.Range("S23").Value = '12/08/2023 08:42'
$ws.Range("T23").Value = 2.26
$ws.Range("U23").Value = '12/08/2023 16:43'
$ws.Range("V23").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/spisska-nova-ves-zilina/prsO5K2C/'

# Row 30
$ws.Range("F30").Value = 'Zilina B'
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 'L. Mikulas'
$ws.Range("I30").Value = 2
$ws.Range("J30").Value = 2.27
$ws.Range("K30").Value = '18/08/2023 02:42'
$ws.Range("L30").Value = 2.29
$ws.Range("M30").Value = '20/08/2023 10:25'
$ws.Range("N30").Value = 3.51
$ws.Range("O30").Value = '18/08/2023 02:42'
$ws.Range("P30").Value = 3.93
$ws.Range("Q30").Value = '20/08/2023 09:55'
$ws.Range("R30").Value = 2.58
$ws.Range("S30").Value = '18/08/2023 02:42'
$ws.Range("T30").Value = 2.63
$ws.Range("U30").Value = '20/08/2023 10:25'
$ws.Range("V30").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/zilina-l-mikulas/vuy6GyHU/'

# Row 31
$ws.Range("F31").Value = 'Slovan Bratislava B'
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 'FK Humenne'
$ws.Range("I31").Value = 2
$ws.Range("J31").Value = 2.69
$ws.Range("K31").Value = '18/08/2023 22:42'
$ws.Range("L31").Value = 2.42
$ws.Range("M31").Value = '20/08/2023 10:28'
$ws.Range("N31").Value = 3.18
$ws.Range("O31").Value = '18/08/2023 22:42'
$ws.Range("P31").Value = 3.44
$ws.Range("Q31").Value = '20/08/2023 10:28'
$ws.Range("R31").Value = 2.34
$ws.Range("S31").Value = '18/08/2023 22:42'
$ws.Range("T31").Value = 2.72
$ws.Range("U31").Value = '20/08/2023 10:28'
$ws.Range("V31").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/slovan-bratislava-fk-humenne/AaqkKJH5/'

# Row 55
$ws.Range("F55").Value = 'Malzenice'
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 'Zilina B'
$ws.Range("I55").Value = 2
$ws.Range("J55").Value = 2.51
$ws.Range("K55").Value = '08/09/2023 22:42'
$ws.Range("L55").Value = 2.51
$ws.Range("M55").Value = '10/09/2023 10:10'
$ws.Range("N55").Value = 3.47
$ws.Range("O55").Value = '08/09/2023 22:42'
$ws.Range("P55").Value = 3.99
$ws.Range("Q55").Value = '10/09/2023 10:10'
$ws.Range("R55").Value = 2.34
$ws.Range("S55").Value = '08/09/2023 22:42'
$ws.Range("T55").Value = 2.37
$ws.Range("U55").Value = '10/09/2023 10:10'
$ws.Range("V55").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/malzenice-zilina/ljKnVemL/'

# Row 56
$ws.Range("F56").Value = 'Petrzalka'
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 'Trebisov'
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 1.49
$ws.Range("K56").Value = '08/09/2023 22:42'
$ws.Range("L56").Value = 1.44
$ws.Range("M56").Value = '10/09/2023 10:14'
$ws.Range("N56").Value = 4.25
$ws.Range("O56").Value = '08/09/2023 22:42'
$ws.Range("P56").Value = 4.7
$ws.Range("Q56").Value = '10/09/2023 10:14'
$ws.Range("R56").Value = 5.03
$ws.Range("S56").Value = '08/09/2023 22:42'
$ws.Range("T56").Value = 6.07
$ws.Range("U56").Value = '10/09/2023 10:14'
$ws.Range("V56").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/petrzalka-trebisov/IPOjUF2R/'

# Row 75
$ws.Range("F75").Value = 'D. Kubin'
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 'Malzenice'
$ws.Range("I75").Value = 3
$ws.Range("J75").Value = 2.02
$ws.Range("K75").Value = '29/09/2023 02:42'
$ws.Range("L75").Value = 2
$ws.Range("M75").Value = '30/09/2023 15:22'
$ws.Range("N75").Value = 3.36
$ws.Range("O75").Value = '29/09/2023 02:42'
$ws.Range("P75").Value = 3.61
$ws.Range("Q75").Value = '30/09/2023 15:22'
$ws.Range("R75").Value = 3.1
$ws.Range("S75").Value = '29/09/2023 02:42'
$ws.Range("T75").Value = 3.39
$ws.Range("U75").Value = '30/09/2023 15:22'
$ws.Range("V75").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-malzenice/tOFzkjvf/'

# Row 76
$ws.Range("F76").Value = 'Puchov'
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 'L. Mikulas'
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 2.16
$ws.Range("K76").Value = '29/09/2023 02:42'
$ws.Range("L76").Value = 2.4
$ws.Range("M76").Value = '30/09/2023 15:28'
$ws.Range("N76").Value = 3.39
$ws.Range("O76").Value = '29/09/2023 02:42'
$ws.Range("P76").Value = 3.45
$ws.Range("Q76").Value = '30/09/2023 15:28'
$ws.Range("R76").Value = 2.81
$ws.Range("S76").Value = '29/09/2023 02:42'
$ws.Range("T76").Value = 2.74
$ws.Range("U76").Value = '30/09/2023 15:28'
$ws.Range("V76").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-l-mikulas/WlQ0eh1Q/'

# Row 77
$ws.Range("F77").Value = 'Samorin'
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 'Komarno'
$ws.Range("I77").Value = 4
$ws.Range("J77").Value = 3.75
$ws.Range("K77").Value = '29/09/2023 21:42'
$ws.Range("L77").Value = 4.36
$ws.Range("M77").Value = '01/10/2023 10:22'
$ws.Range("N77").Value = 3.65
$ws.Range("O77").Value = '29/09/2023 21:42'
$ws.Range("P77").Value = 3.88
$ws.Range("Q77").Value = '01/10/2023 10:22'
$ws.Range("R77").Value = 1.74
$ws.Range("S77").Value = '29/09/2023 21:42'
$ws.Range("T77").Value = 1.71
$ws.Range("U77").Value = '01/10/2023 10:22'
$ws.Range("V77").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/samorin-komarno/On5ulAg0/'

# Row 78
$ws.Range("F78").Value = 'Petrzalka'
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 'FK Humenne'
$ws.Range("I78").Value = 1
$ws.Range("J78").Value = 1.73
$ws.Range("K78").Value = '29/09/2023 21:42'
$ws.Range("L78").Value = 1.65
$ws.Range("M78").Value = '01/10/2023 10:26'
$ws.Range("N78").Value = 3.63
$ws.Range("O78").Value = '29/09/2023 21:42'
$ws.Range("P78").Value = 4.1
$ws.Range("Q78").Value = '01/10/2023 10:28'
$ws.Range("R78").Value = 3.82
$ws.Range("S78").Value = '29/09/2023 21:42'
$ws.Range("T78").Value = 4.53
$ws.Range("U78").Value = '01/10/2023 10:26'
$ws.Range("V78").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/petrzalka-fk-humenne/dEInnlOC/'

# Row 79
$ws.Range("F79").Value = 'Slovan Bratislava B'
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 'Pohronie'
$ws.Range("I79").Value = 3
$ws.Range("J79").Value = 2.23
$ws.Range("K79").Value = '29/09/2023 21:42'
$ws.Range("L79").Value = 2.37
$ws.Range("M79").Value = '01/10/2023 10:24'
$ws.Range("N79").Value = 3.32
$ws.Range("O79").Value = '29/09/2023 21:42'
$ws.Range("P79").Value = 3.63
$ws.Range("Q79").Value = '01/10/2023 10:24'
$ws.Range("R79").Value = 2.74
$ws.Range("S79").Value = '29/09/2023 21:42'
$ws.Range("T79").Value = 2.67
$ws.Range("U79").Value = '01/10/2023 10:24'
$ws.Range("V79").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/slovan-bratislava-fk-pohronie/2w4qmU86/'

# Row 91
$ws.Range("F91").Value = 'Povazska Bystrica'
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 'Zilina B'
$ws.Range("I91").Value = 1
$ws.Range("J91").Value = 1.77
$ws.Range("K91").Value = '13/10/2023 02:13'
$ws.Range("L91").Value = 2.15
$ws.Range("M91").Value = '14/10/2023 14:56'
$ws.Range("N91").Value = 3.81
$ws.Range("O91").Value = '13/10/2023 02:13'
$ws.Range("P91").Value = 3.74
$ws.Range("Q91").Value = '14/10/2023 14:56'
$ws.Range("R91").Value = 3.47
$ws.Range("S91").Value = '13/10/2023 02:13'
$ws.Range("T91").Value = 2.95
$ws.Range("U91").Value = '14/10/2023 14:56'
$ws.Range("V91").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/povazska-bystrica-zilina/dCmt6rFo/'

# Row 92
$ws.Range("F92").Value = 'D. Kubin'
$ws.Range("G92").Value = 2
$ws.Range("H92").Value = 'Spisska Nova Ves'
$ws.Range("I92").Value = 1
$ws.Range("J92").Value = 2
$ws.Range("K92").Value = '13/10/2023 02:13'
$ws.Range("L92").Value = 2.45
$ws.Range("M92").Value = '14/10/2023 14:48'
$ws.Range("N92").Value = 3.36
$ws.Range("O92").Value = '13/10/2023 02:13'
$ws.Range("P92").Value = 3.49
$ws.Range("Q92").Value = '14/10/2023 14:51'
$ws.Range("R92").Value = 3.15
$ws.Range("S92").Value = '13/10/2023 02:13'
$ws.Range("T92").Value = 2.65
$ws.Range("U92").Value = '14/10/2023 14:48'
$ws.Range("V92").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-spisska-nova-ves/IRyk4Mqb/'

# Row 93
$ws.Range("F93").Value = 'FK Humenne'
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 'Malzenice'
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = 1.53
$ws.Range("K93").Value = '13/10/2023 02:13'
$ws.Range("L93").Value = 1.39
$ws.Range("M93").Value = '14/10/2023 14:39'
$ws.Range("N93").Value = 3.95
$ws.Range("O93").Value = '13/10/2023 02:13'
$ws.Range("P93").Value = 4.47
$ws.Range("Q93").Value = '14/10/2023 14:39'
$ws.Range("R93").Value = 4.82
$ws.Range("S93").Value = '13/10/2023 02:13'
$ws.Range("T93").Value = 7.94
$ws.Range("U93").Value = '14/10/2023 14:39'
$ws.Range("V93").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/fk-humenne-malzenice/z3np52Ui/'

# Row 94
$ws.Range("F94").Value = 'Puchov'
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 'Komarno'
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = 2.69
$ws.Range("K94").Value = '13/10/2023 02:13'
$ws.Range("L94").Value = 2.36
$ws.Range("M94").Value = '14/10/2023 14:52'
$ws.Range("N94").Value = 3.23
$ws.Range("O94").Value = '13/10/2023 02:13'
$ws.Range("P94").Value = 3.38
$ws.Range("Q94").Value = '14/10/2023 14:52'
$ws.Range("R94").Value = 2.31
$ws.Range("S94").Value = '13/10/2023 02:13'
$ws.Range("T94").Value = 2.84
$ws.Range("U94").Value = '14/10/2023 14:52'
$ws.Range("V94").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-komarno/8dfbsaxo/'

# Row 98
$ws.Range("F98").Value = 'D. Kubin'
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 'L. Mikulas'
$ws.Range("I98").Value = 2
$ws.Range("J98").Value = 6.42
$ws.Range("K98").Value = '20/10/2023 01:53'
$ws.Range("L98").Value = 4.12
$ws.Range("M98").Value = '21/10/2023 14:29'
$ws.Range("N98").Value = 4.84
$ws.Range("O98").Value = '20/10/2023 01:53'
$ws.Range("P98").Value = 3.58
$ws.Range("Q98").Value = '21/10/2023 14:29'
$ws.Range("R98").Value = 1.33
$ws.Range("S98").Value = '20/10/2023 01:53'
$ws.Range("T98").Value = 1.82
$ws.Range("U98").Value = '21/10/2023 14:29'
$ws.Range("V98").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-l-mikulas/tQtUfI6p/'

# Row 99
$ws.Range("F99").Value = 'Komarno'
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 'Trebisov'
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1.18
$ws.Range("K99").Value = '20/10/2023 01:42'
$ws.Range("L99").Value = 1.28
$ws.Range("M99").Value = '21/10/2023 14:15'
$ws.Range("N99").Value = 6.03
$ws.Range("O99").Value = '20/10/2023 01:42'
$ws.Range("P99").Value = 5.64
$ws.Range("Q99").Value = '21/10/2023 14:17'
$ws.Range("R99").Value = 9.76
$ws.Range("S99").Value = '20/10/2023 01:42'
$ws.Range("T99").Value = 9.09
$ws.Range("U99").Value = '21/10/2023 14:15'
$ws.Range("V99").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/komarno-trebisov/pfeb0vqN/'

# Row 100
$ws.Range("F100").Value = 'Malzenice'
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = 'Myjava'
$ws.Range("I100").Value = 1
$ws.Range("J100").Value = 2.59
$ws.Range("K100").Value = '20/10/2023 01:42'
$ws.Range("L100").Value = 3.04
$ws.Range("M100").Value = '21/10/2023 14:28'
$ws.Range("N100").Value = 3.23
$ws.Range("O100").Value = '20/10/2023 01:42'
$ws.Range("P100").Value = 3.14
$ws.Range("Q100").Value = '21/10/2023 14:28'
$ws.Range("R100").Value = 2.4
$ws.Range("S100").Value = '20/10/2023 01:42'
$ws.Range("T100").Value = 2.36
$ws.Range("U100").Value = '21/10/2023 14:23'
$ws.Range("V100").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/malzenice-myjava/4bSOGaT3/'

# Row 101
$ws.Range("F101").Value = 'Presov'
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 'Petrzalka'
$ws.Range("I101").Value = 1
$ws.Range("J101").Value = 1.56
$ws.Range("K101").Value = '20/10/2023 01:42'
$ws.Range("L101").Value = 1.85
$ws.Range("M101").Value = '21/10/2023 14:19'
$ws.Range("N101").Value = 3.84
$ws.Range("O101").Value = '20/10/2023 01:42'
$ws.Range("P101").Value = 3.74
$ws.Range("Q101").Value = '21/10/2023 14:19'
$ws.Range("R101").Value = 4.69
$ws.Range("S101").Value = '20/10/2023 01:42'
$ws.Range("T101").Value = 3.82
$ws.Range("U101").Value = '21/10/2023 14:19'
$ws.Range("V101").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/presov-petrzalka/Opf2abbT/'

# Row 102
$ws.Range("F102").Value = 'Spisska Nova Ves'
$ws.Range("G102").Value = 2
$ws.Range("H102").Value = 'FK Humenne'
$ws.Range("I102").Value = 3
$ws.Range("J102").Value = 3.02
$ws.Range("K102").Value = '20/10/2023 01:42'
$ws.Range("L102").Value = 3.43
$ws.Range("M102").Value = '21/10/2023 14:27'
$ws.Range("N102").Value = 3.19
$ws.Range("O102").Value = '20/10/2023 01:42'
$ws.Range("P102").Value = 3.47
$ws.Range("Q102").Value = '21/10/2023 14:27'
$ws.Range("R102").Value = 2.13
$ws.Range("S102").Value = '20/10/2023 01:42'
$ws.Range("T102").Value = 2.04
$ws.Range("U102").Value = '21/10/2023 14:27'
$ws.Range("V102").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/spisska-nova-ves-fk-humenne/fTpYgxMj/'

# Row 105
$ws.Range("F105").Value = 'Povazska Bystrica'
$ws.Range("G105").Value = 3
$ws.Range("H105").Value = 'Slovan Bratislava B'
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 1.8
$ws.Range("K105").Value = '27/10/2023 02:42'
$ws.Range("L105").Value = 1.46
$ws.Range("M105").Value = '28/10/2023 13:57'
$ws.Range("N105").Value = 3.56
$ws.Range("O105").Value = '27/10/2023 02:42'
$ws.Range("P105").Value = 4.44
$ws.Range("Q105").Value = '28/10/2023 13:57'
$ws.Range("R105").Value = 3.58
$ws.Range("S105").Value = '27/10/2023 02:42'
$ws.Range("T105").Value = 6.15
$ws.Range("U105").Value = '28/10/2023 13:57'
$ws.Range("V105").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/povazska-bystrica-slovan-bratislava/ObZzEcDM/'

# Row 106
$ws.Range("F106").Value = 'Puchov'
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 'Pohronie'
$ws.Range("I106").Value = 4
$ws.Range("J106").Value = 1.75
$ws.Range("K106").Value = '27/10/2023 02:42'
$ws.Range("L106").Value = 1.71
$ws.Range("M106").Value = '28/10/2023 14:21'
$ws.Range("N106").Value = 3.6
$ws.Range("O106").Value = '27/10/2023 02:42'
$ws.Range("P106").Value = 3.81
$ws.Range("Q106").Value = '28/10/2023 14:21'
$ws.Range("R106").Value = 3.77
$ws.Range("S106").Value = '27/10/2023 02:42'
$ws.Range("T106").Value = 4.45
$ws.Range("U106").Value = '28/10/2023 14:21'
$ws.Range("V106").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-fk-pohronie/QDGXEwbG/'

# Row 107
$ws.Range("F107").Value = 'D. Kubin'
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 'Presov'
$ws.Range("I107").Value = 3
$ws.Range("J107").Value = 5.32
$ws.Range("K107").Value = '27/10/2023 02:42'
$ws.Range("L107").Value = 8.949999999999999
$ws.Range("M107").Value = '28/10/2023 14:28'
$ws.Range("N107").Value = 3.95
$ws.Range("O107").Value = '27/10/2023 02:42'
$ws.Range("P107").Value = 5.83
$ws.Range("Q107").Value = '28/10/2023 14:28'
$ws.Range("R107").Value = 1.48
$ws.Range("S107").Value = '27/10/2023 02:42'
$ws.Range("T107").Value = 1.26
$ws.Range("U107").Value = '28/10/2023 14:28'
$ws.Range("V107").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-presov/ILYvDHSS/'

# Row 108
$ws.Range("F108").Value = 'FK Humenne'
$ws.Range("G108").Value = 4
$ws.Range("H108").Value = 'L. Mikulas'
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = 2.3
$ws.Range("K108").Value = '27/10/2023 02:42'
$ws.Range("L108").Value = 1.96
$ws.Range("M108").Value = '28/10/2023 14:21'
$ws.Range("N108").Value = 3.18
$ws.Range("O108").Value = '27/10/2023 02:42'
$ws.Range("P108").Value = 3.59
$ws.Range("Q108").Value = '28/10/2023 14:21'
$ws.Range("R108").Value = 2.74
$ws.Range("S108").Value = '27/10/2023 02:42'
$ws.Range("T108").Value = 3.52
$ws.Range("U108").Value = '28/10/2023 14:21'
$ws.Range("V108").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/fk-humenne-l-mikulas/x6FB7ejj/'

# ---- Append new rows 113-116 ----
# Copy formatting (styles) from the last existing data row (112)
$ws.Range("A112:V112").Copy()
$ws.Range("A113:V116").PasteSpecial(-4122)

# Row 113
$ws.Range("A113").Value = 112
$ws.Range("B113").Value = 'slovakia'
$ws.Range("C113").Value = '2-liga'
$ws.Range("D113").Value = '2023-2024'
$ws.Range("E113").Value = 45234.54166666666
$ws.Range("F113").Value = 'Spisska Nova Ves'
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 'L. Mikulas'
$ws.Range("I113").Value = 2
$ws.Range("J113").Value = 3.71
$ws.Range("K113").Value = '03/11/2023 01:13'
$ws.Range("L113").Value = 4.49
$ws.Range("M113").Value = '04/11/2023 12:50'
$ws.Range("N113").Value = 3.48
$ws.Range("O113").Value = '03/11/2023 01:13'
$ws.Range("P113").Value = 3.68
$ws.Range("Q113").Value = '04/11/2023 12:50'
$ws.Range("R113").Value = 1.79
$ws.Range("S113").Value = '03/11/2023 01:13'
$ws.Range("T113").Value = 1.73
$ws.Range("U113").Value = '04/11/2023 12:50'
$ws.Range("V113").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/spisska-nova-ves-l-mikulas/jZhIQzcS/'

# Row 114
$ws.Range("A114").Value = 113
$ws.Range("B114").Value = 'slovakia'
$ws.Range("C114").Value = '2-liga'
$ws.Range("D114").Value = '2023-2024'
$ws.Range("E114").Value = 45234.54166666666
$ws.Range("F114").Value = 'Komarno'
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 'D. Kubin'
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 1.17
$ws.Range("K114").Value = '03/11/2023 01:12'
$ws.Range("L114").Value = 1.21
$ws.Range("M114").Value = '04/11/2023 12:49'
$ws.Range("N114").Value = 6.27
$ws.Range("O114").Value = '03/11/2023 01:12'
$ws.Range("P114").Value = 6.5
$ws.Range("Q114").Value = '04/11/2023 12:54'
$ws.Range("R114").Value = 9.890000000000001
$ws.Range("S114").Value = '03/11/2023 01:12'
$ws.Range("T114").Value = 10.84
$ws.Range("U114").Value = '04/11/2023 12:54'
$ws.Range("V114").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/komarno-d-kubin/4IfASdSF/'

# Row 115
$ws.Range("A115").Value = 114
$ws.Range("B115").Value = 'slovakia'
$ws.Range("C115").Value = '2-liga'
$ws.Range("D115").Value = '2023-2024'
$ws.Range("E115").Value = 45234.54166666666
$ws.Range("F115").Value = 'Presov'
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 'FK Humenne'
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 1.53
$ws.Range("K115").Value = '03/11/2023 01:12'
$ws.Range("L115").Value = 1.57
$ws.Range("M115").Value = '04/11/2023 12:50'
$ws.Range("N115").Value = 3.81
$ws.Range("O115").Value = '03/11/2023 01:12'
$ws.Range("P115").Value = 3.94
$ws.Range("Q115").Value = '04/11/2023 12:50'
$ws.Range("R115").Value = 5.01
$ws.Range("S115").Value = '03/11/2023 01:12'
$ws.Range("T115").Value = 5.53
$ws.Range("U115").Value = '04/11/2023 12:50'
$ws.Range("V115").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/presov-fk-humenne/CrfERGsM/'

# Row 116
$ws.Range("A116").Value = 115
$ws.Range("B116").Value = 'slovakia'
$ws.Range("C116").Value = '2-liga'
$ws.Range("D116").Value = '2023-2024'
$ws.Range("E116").Value = 45234.60416666666
$ws.Range("F116").Value = 'Zilina B'
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 'Petrzalka'
$ws.Range("I116").Value = 1
$ws.Range("J116").Value = 2.55
$ws.Range("K116").Value = '03/11/2023 02:42'
$ws.Range("L116").Value = 2.55
$ws.Range("M116").Value = '04/11/2023 14:21'
$ws.Range("N116").Value = 3.41
$ws.Range("O116").Value = '03/11/2023 02:42'
$ws.Range("P116").Value = 3.86
$ws.Range("Q116").Value = '04/11/2023 14:29'
$ws.Range("R116").Value = 2.39
$ws.Range("S116").Value = '03/11/2023 02:42'
$ws.Range("T116").Value = 2.38
$ws.Range("U116").Value = '04/11/2023 14:29'
$ws.Range("V116").Value = 'https://www.betexplorer.com/football/slovakia/2-liga/zilina-petrzalka/dGVpMf5q/'
